$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 61: fill in the previously-missing base mod columns (C, D, E) ---
# (F, G, H, I, J already contained 0.34 for this row)
$ws.Range("C61").Value = 0.34
$ws.Range("D61").Value = 0.34
$ws.Range("E61").Value = 0.34

# --- New affix rows 76-90 ---

# Row 76: Earth Tuned (suffix)
$ws.Range("A76").Value = "Earth Tuned"
$ws.Range("B76").Value = "Be tuned into the earth and its balancing energy"
$ws.Range("C76").Value = 0.38
$ws.Range("D76").Value = 0.38
$ws.Range("E76").Value = 0.38
$ws.Range("F76").Value = 0.38
$ws.Range("G76").Value = 0.38
$ws.Range("H76").Value = 0.38
$ws.Range("I76").Value = 0.38
$ws.Range("J76").Value = 0.38
$ws.Range("K76").Value = 500
$ws.Range("L76").Value = 110
$ws.Range("M76").Value = 200
$ws.Range("Q76").Value = 2600980000
$ws.Range("R76").Value = "suffix"

# Row 77: Strength of Courage (suffix)
$ws.Range("A77").Value = "Strength of Courage"
$ws.Range("B77").Value = "Give your self the belief in strength and courage to survive whats to come."
$ws.Range("F77").Value = 0.48
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 530
$ws.Range("L77").Value = 120
$ws.Range("M77").Value = 200
$ws.Range("Q77").Value = 5670200000
$ws.Range("R77").Value = "suffix"

# Row 78: See All (suffix)
$ws.Range("A78").Value = "See All"
$ws.Range("B78").Value = "Missing never was the option. Damage and accuracy was."
$ws.Range("H78").Value = 0.48
$ws.Range("K78").Value = 500
$ws.Range("L78").Value = 120
$ws.Range("M78").Value = 200
$ws.Range("Q78").Value = 5670200000
$ws.Range("R78").Value = "suffix"

# Row 79: Reapers Kiss (suffix)
$ws.Range("A79").Value = "Reapers Kiss"
$ws.Range("B79").Value = "Thirst for souls, thirst for death."
$ws.Range("G79").Value = 0.48
$ws.Range("K79").Value = 500
$ws.Range("L79").Value = 120
$ws.Range("M79").Value = 200
$ws.Range("Q79").Value = 5670200000
$ws.Range("R79").Value = "suffix"

# Row 80: Scholars Research (suffix)
$ws.Range("A80").Value = "Scholars Research"
$ws.Range("B80").Value = "All the research in the world that has been done so far on the study of magic, is yours."
$ws.Range("J80").Value = 0.48
$ws.Range("K80").Value = 500
$ws.Range("L80").Value = 120
$ws.Range("M80").Value = 200
$ws.Range("Q80").Value = 5670200000
$ws.Range("R80").Value = "suffix"

# Row 81: Prayer Of Hope (suffix)
$ws.Range("A81").Value = "Prayer Of Hope"
$ws.Range("B81").Value = "If you pray, I shall give you hope. If you listen, I shall guide you."
$ws.Range("I81").Value = 0.48
$ws.Range("K81").Value = 500
$ws.Range("L81").Value = 120
$ws.Range("M81").Value = 200
$ws.Range("Q81").Value = 5670200000
$ws.Range("R81").Value = "suffix"

# Row 82: Rangers Luck (prefix, Accuracy)
$ws.Range("A82").Value = "Rangers Luck"
$ws.Range("B82").Value = "With the luck of a well trained ranger you might be able to finally hit that bear."
$ws.Range("H82").Value = 0.28
$ws.Range("K82").Value = 600
$ws.Range("L82").Value = 130
$ws.Range("M82").Value = 200
$ws.Range("N82").Value = "Accuracy"
$ws.Range("O82").Value = 0.15
$ws.Range("P82").Value = 0.18
$ws.Range("Q82").Value = 7568900000
$ws.Range("R82").Value = "prefix"

# Row 83: Rumor's Movement (prefix, Dodge)
$ws.Range("A83").Value = "Rumor's Movement"
$ws.Range("B83").Value = "Move like a Rumor. Impossible o see, hit or know is coming."
$ws.Range("H83").Value = 0.28
$ws.Range("K83").Value = 600
$ws.Range("L83").Value = 130
$ws.Range("M83").Value = 200
$ws.Range("N83").Value = "Dodge"
$ws.Range("O83").Value = 0.15
$ws.Range("P83").Value = 0.18
$ws.Range("Q83").Value = 7568900000
$ws.Range("R83").Value = "prefix"

# Row 84: Smell of Gold (prefix, Looting)
$ws.Range("A84").Value = "Smell of Gold"
$ws.Range("B84").Value = "The smell of gold can drag one deep into the lust for treasure."
$ws.Range("K84").Value = 600
$ws.Range("L84").Value = 130
$ws.Range("M84").Value = 200
$ws.Range("N84").Value = "Looting"
$ws.Range("O84").Value = 0.15
$ws.Range("P84").Value = 0.18
$ws.Range("Q84").Value = 7568900000
$ws.Range("R84").Value = "prefix"

# Row 85: Godly Weapon Crafting (prefix, Weapon Crafting)
$ws.Range("A85").Value = "Godly Weapon Crafting"
$ws.Range("B85").Value = "You'll probably never fail to craft again."
$ws.Range("K85").Value = 700
$ws.Range("L85").Value = 140
$ws.Range("M85").Value = 200
$ws.Range("N85").Value = "Weapon Crafting"
$ws.Range("O85").Value = 0.4
$ws.Range("P85").Value = 0.45
$ws.Range("Q85").Value = 10670900000
$ws.Range("R85").Value = "prefix"

# Row 86: Angelic Armour Smithing (prefix, Armour Crafting)
$ws.Range("A86").Value = "Angelic Armour Smithing"
$ws.Range("B86").Value = "Craft amour with the help of the angels."
$ws.Range("K86").Value = 700
$ws.Range("L86").Value = 140
$ws.Range("M86").Value = 200
$ws.Range("N86").Value = "Armour Crafting"
$ws.Range("O86").Value = 0.4
$ws.Range("P86").Value = 0.45
$ws.Range("Q86").Value = 10670900000
$ws.Range("R86").Value = "prefix"

# Row 87: Divine Magic Crafting (prefix, Spell Crafting)
$ws.Range("A87").Value = "Divine Magic Crafting"
$ws.Range("B87").Value = "Use divine help to craft magical spells."
$ws.Range("K87").Value = 700
$ws.Range("L87").Value = 140
$ws.Range("M87").Value = 200
$ws.Range("N87").Value = "Spell Crafting"
$ws.Range("O87").Value = 0.4
$ws.Range("P87").Value = 0.45
$ws.Range("Q87").Value = 10670900000
$ws.Range("R87").Value = "prefix"

# Row 88: Devilish Ring Crafting (prefix, Ring Crafting)
$ws.Range("A88").Value = "Devilish Ring Crafting"
$ws.Range("B88").Value = "Let the thoughts and the inspiration from the devil inspire you."
$ws.Range("K88").Value = 700
$ws.Range("L88").Value = 140
$ws.Range("M88").Value = 200
$ws.Range("N88").Value = "Ring Crafting"
$ws.Range("O88").Value = 0.4
$ws.Range("P88").Value = 0.45
$ws.Range("Q88").Value = 10670900000
$ws.Range("R88").Value = "prefix"

# Row 89: Enchanted Labyrinth (prefix, Enchanting)
$ws.Range("A89").Value = "Enchanted Labyrinth"
$ws.Range("B89").Value = "Get lost in a labyrinth of all the enchantments you can create."
$ws.Range("K89").Value = 700
$ws.Range("L89").Value = 120
$ws.Range("M89").Value = 200
$ws.Range("N89").Value = "Enchanting"
$ws.Range("O89").Value = 0.4
$ws.Range("P89").Value = 0.45
$ws.Range("Q89").Value = 10670900000
$ws.Range("R89").Value = "prefix"

# Row 90: Dark Pact (prefix, Artifact Crafting)
$ws.Range("A90").Value = "Dark Pact"
$ws.Range("B90").Value = "Make a dark pact when creating artifacts."
$ws.Range("K90").Value = 700
$ws.Range("L90").Value = 120
$ws.Range("M90").Value = 200
$ws.Range("N90").Value = "Artifact Crafting"
$ws.Range("O90").Value = 0.4
$ws.Range("P90").Value = 0.45
$ws.Range("Q90").Value = 10670900000
$ws.Range("R90").Value = "prefix"

# --- Column Q got a little wider to fit the new, larger cost values ---
$ws.Range("Q1").EntireColumn.ColumnWidth = 12.1
